# Regenerate the "K" column (column G) values for rows 2-31 on the active
# worksheet, matching the regenerated save_data (K computed from pitch-level
# Strike# data instead of the prior approximation).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 2
    4  = 0
    5  = 2
    6  = 2
    7  = 1
    8  = 1
    9  = 1
    10 = 2
    11 = 0
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 0
    22 = 1
    23 = 4
    24 = 1
    25 = 1
    26 = 3
    27 = 1
    28 = 1
    29 = 3
    30 = 2
    31 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
